$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp text in A1
$ws.Range("A1").Value = "Datos actualizados a 4 de Abril de 2020 a las 04:20"

# Swap Honduras in ahead of Oman (rows 95-99): Honduras new stats, others shift down one slot
$ws.Range("A95").Value = "Honduras"
$ws.Cells.Item(95,2).Value = 264
$ws.Cells.Item(95,3).Value = 42
$ws.Cells.Item(95,4).Value = 3
$ws.Cells.Item(95,5).Value = 246
$ws.Cells.Item(95,6).Value = 10
$ws.Cells.Item(95,7).Value = 0
$ws.Cells.Item(95,8).Value = 15

$ws.Range("A96").Value = "Oman"
$ws.Cells.Item(96,2).Value = 252
$ws.Cells.Item(96,3).Value = 0
$ws.Cells.Item(96,4).Value = 57
$ws.Cells.Item(96,5).Value = 194
$ws.Cells.Item(96,6).Value = 3
$ws.Cells.Item(96,7).Value = 0
$ws.Cells.Item(96,8).Value = 1

$ws.Range("A97").Value = "San Marino"
$ws.Cells.Item(97,2).Value = 251
$ws.Cells.Item(97,3).Value = 0
$ws.Cells.Item(97,4).Value = 26
$ws.Cells.Item(97,5).Value = 193
$ws.Cells.Item(97,6).Value = 13
$ws.Cells.Item(97,7).Value = 0
$ws.Cells.Item(97,8).Value = 32

$ws.Range("A98").Value = "Vietnam"
$ws.Cells.Item(98,2).Value = 239
$ws.Cells.Item(98,3).Value = 0
$ws.Cells.Item(98,4).Value = 85
$ws.Cells.Item(98,5).Value = 154
$ws.Cells.Item(98,6).Value = 3
$ws.Cells.Item(98,7).Value = 0
$ws.Cells.Item(98,8).Value = 0

$ws.Range("A99").Value = "Uzbekistan"
$ws.Cells.Item(99,2).Value = 227
$ws.Cells.Item(99,3).Value = 0
$ws.Cells.Item(99,4).Value = 25
$ws.Cells.Item(99,5).Value = 200
$ws.Cells.Item(99,6).Value = 8
$ws.Cells.Item(99,7).Value = 0
$ws.Cells.Item(99,8).Value = 2

# Sri Lanka (row 109) stats update only
$ws.Cells.Item(109,2).Value = 159
$ws.Cells.Item(109,3).Value = 0
$ws.Cells.Item(109,4).Value = 24
$ws.Cells.Item(109,5).Value = 130
$ws.Cells.Item(109,6).Value = 5
$ws.Cells.Item(109,7).Value = 1
$ws.Cells.Item(109,8).Value = 5

# Swap Paraguay in ahead of Gibraltar (rows 123-125)
$ws.Range("A123").Value = "Trinidad yTobago"
$ws.Cells.Item(123,2).Value = 100
$ws.Cells.Item(123,3).Value = 2
$ws.Cells.Item(123,4).Value = 1
$ws.Cells.Item(123,5).Value = 93
$ws.Cells.Item(123,6).Value = 0
$ws.Cells.Item(123,7).Value = 0
$ws.Cells.Item(123,8).Value = 6

$ws.Range("A124").Value = "Paraguay"
$ws.Cells.Item(124,2).Value = 96
$ws.Cells.Item(124,3).Value = 4
$ws.Cells.Item(124,4).Value = 12
$ws.Cells.Item(124,5).Value = 81
$ws.Cells.Item(124,6).Value = 2
$ws.Cells.Item(124,7).Value = 0
$ws.Cells.Item(124,8).Value = 3

$ws.Range("A125").Value = "Gibraltar"
$ws.Cells.Item(125,2).Value = 95
$ws.Cells.Item(125,3).Value = 0
$ws.Cells.Item(125,4).Value = 46
$ws.Cells.Item(125,5).Value = 49
$ws.Cells.Item(125,6).Value = 0
$ws.Cells.Item(125,7).Value = 0
$ws.Cells.Item(125,8).Value = 0
